# Support sheet names with spaces:
#   "data"  -> "data with charts"
#   "data2" -> "charts"
# Renaming the sheets automatically updates every worksheet formula and
# defined name that referred to them (Excel quotes sheet names that
# contain spaces, e.g. 'data with charts'!$A$1).

$wb = $excel.ActiveWorkbook

$oldDataName = "data"
$newDataName = "data with charts"
$oldChartsName = "data2"
$newChartsName = "charts"

$wb.Worksheets.Item($oldDataName).Name = $newDataName
$wb.Worksheets.Item($oldChartsName).Name = $newChartsName

# Sheet renames don't retroactively rewrite the cached formula text that
# charts keep in their series (c:f elements) - those still point at the
# old, unquoted "data!..." references. Walk every chart on every sheet
# and rewrite any series formula that still mentions the old name so it
# points at the new, quoted sheet name instead.
foreach ($sheet in $wb.Worksheets) {
    $chartObjects = $sheet.ChartObjects()
    for ($i = 1; $i -le $chartObjects.Count; $i++) {
        $chart = $chartObjects.Item($i).Chart
        $series = $chart.SeriesCollection()
        for ($j = 1; $j -le $series.Count; $j++) {
            $ser = $series.Item($j)
            $formula = $ser.Formula
            $updated = $formula.Replace("$oldDataName!", "'$newDataName'!")
            $updated = $updated.Replace("$oldChartsName!", "$newChartsName!")
            if ($updated -ne $formula) {
                $ser.Formula = $updated
            }
        }
    }
}

# The "charts" sheet (formerly "data2") becomes the active tab, with a
# new selection on it; the "data with charts" sheet (formerly "data")
# loses its tabSelected flag, keeping its existing selection.
$chartsSheet = $wb.Worksheets.Item($newChartsName)
$chartsSheet.Activate()
$chartsSheet.Range("E14").Select()
